$d = $word.ActiveDocument

$xml0 = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="203E6B75" w14:textId="06D67B88" w:rsidR="001D7888" w:rsidRPr="00B87C48" w:rsidRDefault="001D7888" w:rsidP="001D7888"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="13"/></w:numPr><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Dax-Regular" w:hAnsi="Dax-Regular"/><w:lang w:val="es-419"/></w:rPr></w:pPr><w:r w:rsidRPr="001D7888"><w:rPr><w:lang w:val="es-419"/></w:rPr><w:t xml:space="preserve">¿Qué configuración de ideal ADT Map escogería para el </w:t></w:r><w:r w:rsidRPr="001D7888"><w:rPr><w:b/><w:bCs/><w:lang w:val="es-419"/></w:rPr><w:t xml:space="preserve">índice </w:t></w:r><w:r w:rsidR="001978BE"><w:rPr><w:b/><w:bCs/><w:lang w:val="es-419"/></w:rPr><w:t>de años (“Año”)</w:t></w:r><w:r w:rsidR="001978BE" w:rsidRPr="001D7888"><w:rPr><w:lang w:val="es-419"/></w:rPr><w:t xml:space="preserve"> ?,</w:t></w:r><w:r w:rsidRPr="001D7888"><w:rPr><w:lang w:val="es-419"/></w:rPr><w:t xml:space="preserve"> especifique el mecanismo de colisión, el factor de carga y el numero inicial de elementos.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Paragraphs(205).Range.InsertXML($xml0)

$xml1 = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="3281F67E" w14:textId="77777777" w:rsidR="001D7888" w:rsidRDefault="001D7888" w:rsidP="001D7888"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="13"/></w:numPr><w:jc w:val="both"/><w:rPr><w:lang w:val="es-419"/></w:rPr></w:pPr><w:r w:rsidRPr="001D7888"><w:rPr><w:lang w:val="es-419"/></w:rPr><w:t xml:space="preserve">¿Qué cambios percibe en el </w:t></w:r><w:r w:rsidRPr="001D7888"><w:rPr><w:b/><w:bCs/><w:lang w:val="es-419"/></w:rPr><w:t>consumo de memoria</w:t></w:r><w:r w:rsidRPr="001D7888"><w:rPr><w:lang w:val="es-419"/></w:rPr><w:t xml:space="preserve"> al modificar el esquema de colisiones?, si los percibe, describa las diferencias y argumente su respuesta.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Paragraphs(203).Range.InsertXML($xml1)

$xml2 = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="25E04868" w14:textId="77777777" w:rsidR="001D7888" w:rsidRDefault="001D7888" w:rsidP="001D7888"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="13"/></w:numPr><w:jc w:val="both"/><w:rPr><w:lang w:val="es-419"/></w:rPr></w:pPr><w:r w:rsidRPr="001D7888"><w:rPr><w:lang w:val="es-419"/></w:rPr><w:t xml:space="preserve">¿Qué cambios percibe en el </w:t></w:r><w:r w:rsidRPr="001D7888"><w:rPr><w:b/><w:bCs/><w:lang w:val="es-419"/></w:rPr><w:t>consumo de memoria</w:t></w:r><w:r w:rsidRPr="001D7888"><w:rPr><w:lang w:val="es-419"/></w:rPr><w:t xml:space="preserve"> al modificar el factor de carga máximo para cargar el catálogo de contenido Streaming?</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Paragraphs(199).Range.InsertXML($xml2)

$xml3 = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="219D2D51" w14:textId="77777777" w:rsidR="001D7888" w:rsidRDefault="001D7888" w:rsidP="001D7888"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="13"/></w:numPr><w:jc w:val="both"/><w:rPr><w:lang w:val="es-419"/></w:rPr></w:pPr><w:r w:rsidRPr="001D7888"><w:rPr><w:lang w:val="es-419"/></w:rPr><w:t xml:space="preserve">¿Qué cambios percibe en el </w:t></w:r><w:r w:rsidRPr="001D7888"><w:rPr><w:b/><w:bCs/><w:lang w:val="es-419"/></w:rPr><w:t xml:space="preserve">tiempo de ejecución </w:t></w:r><w:r w:rsidRPr="001D7888"><w:rPr><w:lang w:val="es-419"/></w:rPr><w:t>al modificar el factor de carga máximo para cargar el catálogo de contenido Streaming?</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Paragraphs(197).Range.InsertXML($xml3)

$xml4 = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="5BF6AECF" w14:textId="77777777" w:rsidR="001D7888" w:rsidRDefault="001D7888" w:rsidP="001D7888"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="13"/></w:numPr><w:jc w:val="both"/><w:rPr><w:lang w:val="es-419"/></w:rPr></w:pPr><w:r w:rsidRPr="001D7888"><w:rPr><w:lang w:val="es-419"/></w:rPr><w:t>Dado el número de elementos de los archivos del reto (large), ¿Cuál sería el factor de carga para estos índices según su mecanismo de colisión?</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Paragraphs(195).Range.InsertXML($xml4)

$xml5 = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="0C8B6D8F" w14:textId="77777777" w:rsidR="001D7888" w:rsidRDefault="001D7888" w:rsidP="001D7888"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="13"/></w:numPr><w:jc w:val="both"/><w:rPr><w:lang w:val="es-419"/></w:rPr></w:pPr><w:r w:rsidRPr="001D7888"><w:rPr><w:lang w:val="es-419"/></w:rPr><w:t xml:space="preserve">Según los índices propuestos ¿en qué caso usaría </w:t></w:r><w:r w:rsidRPr="001D7888"><w:rPr><w:b/><w:bCs/><w:lang w:val="es-419"/></w:rPr><w:t>Linear Probing</w:t></w:r><w:r w:rsidRPr="001D7888"><w:rPr><w:lang w:val="es-419"/></w:rPr><w:t xml:space="preserve"> o </w:t></w:r><w:r w:rsidRPr="001D7888"><w:rPr><w:b/><w:bCs/><w:lang w:val="es-419"/></w:rPr><w:t>Separate Chaining</w:t></w:r><w:r w:rsidRPr="001D7888"><w:rPr><w:lang w:val="es-419"/></w:rPr><w:t xml:space="preserve"> en estos índices? y ¿Por qué?</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Paragraphs(193).Range.InsertXML($xml5)

$xml6 = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="7FDD4806" w14:textId="303E6E74" w:rsidR="001D7888" w:rsidRPr="001D7888" w:rsidRDefault="00926514" w:rsidP="00926514"><w:pPr><w:rPr><w:lang w:val="es-419"/></w:rPr></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Paragraphs(190).Range.InsertXML($xml6)

$xml7 = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="4452CB3C" w14:textId="4CF53129" w:rsidR="001D7888" w:rsidRPr="001D7888" w:rsidRDefault="001C62B2" w:rsidP="001C62B2"><w:pPr><w:rPr><w:lang w:val="es-419"/></w:rPr></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Paragraphs(188).Range.InsertXML($xml7)

$xml8 = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="42530D93" w14:textId="77777777" w:rsidR="001D7888" w:rsidRDefault="001D7888" w:rsidP="001D7888"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="13"/></w:numPr><w:jc w:val="both"/><w:rPr><w:lang w:val="es-419"/></w:rPr></w:pPr><w:r w:rsidRPr="001D7888"><w:rPr><w:lang w:val="es-419"/></w:rPr><w:t xml:space="preserve">¿Por qué son importantes las funciones </w:t></w:r><w:r w:rsidRPr="001D7888"><w:rPr><w:b/><w:bCs/><w:lang w:val="es-419"/></w:rPr><w:t>start()</w:t></w:r><w:r w:rsidRPr="001D7888"><w:rPr><w:lang w:val="es-419"/></w:rPr><w:t xml:space="preserve"> y </w:t></w:r><w:r w:rsidRPr="001D7888"><w:rPr><w:b/><w:bCs/><w:lang w:val="es-419"/></w:rPr><w:t>stop()</w:t></w:r><w:r w:rsidRPr="001D7888"><w:rPr><w:lang w:val="es-419"/></w:rPr><w:t xml:space="preserve"> de la librería </w:t></w:r><w:r w:rsidRPr="001D7888"><w:rPr><w:b/><w:bCs/><w:lang w:val="es-419"/></w:rPr><w:t>tracemalloc</w:t></w:r><w:r w:rsidRPr="001D7888"><w:rPr><w:lang w:val="es-419"/></w:rPr><w:t>?</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Paragraphs(187).Range.InsertXML($xml8)

$xml9 = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="50C609D3" w14:textId="409B32F5" w:rsidR="001D7888" w:rsidRPr="001D7888" w:rsidRDefault="001C62B2" w:rsidP="001C62B2"><w:pPr><w:jc w:val="both"/><w:rPr><w:lang w:val="es-419"/></w:rPr></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Paragraphs(186).Range.InsertXML($xml9)

$xml10 = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="49065933" w14:textId="5F88E477" w:rsidR="001D7888" w:rsidRDefault="001D7888" w:rsidP="001D7888"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="13"/></w:numPr><w:jc w:val="both"/><w:rPr><w:lang w:val="es-419"/></w:rPr></w:pPr><w:r w:rsidRPr="001D7888"><w:rPr><w:lang w:val="es-419"/></w:rPr><w:t xml:space="preserve">¿Por qué en la función </w:t></w:r><w:proofErr w:type="gramStart"/><w:r w:rsidRPr="001D7888"><w:rPr><w:b/><w:bCs/><w:lang w:val="es-419"/></w:rPr><w:t>getTime(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r w:rsidRPr="001D7888"><w:rPr><w:b/><w:bCs/><w:lang w:val="es-419"/></w:rPr><w:t>)</w:t></w:r><w:r w:rsidRPr="001D7888"><w:rPr><w:lang w:val="es-419"/></w:rPr><w:t xml:space="preserve"> se utiliza </w:t></w:r><w:r w:rsidRPr="001D7888"><w:rPr><w:b/><w:bCs/><w:lang w:val="es-419"/></w:rPr><w:t>time.perf_counter()</w:t></w:r><w:r w:rsidRPr="001D7888"><w:rPr><w:lang w:val="es-419"/></w:rPr><w:t xml:space="preserve"> en vez de otras funciones como </w:t></w:r><w:r w:rsidRPr="001D7888"><w:rPr><w:b/><w:bCs/><w:lang w:val="es-419"/></w:rPr><w:t>time.process_time()</w:t></w:r><w:r w:rsidRPr="001D7888"><w:rPr><w:lang w:val="es-419"/></w:rPr><w:t>?</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Paragraphs(185).Range.InsertXML($xml10)

$xml11 = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="733D2572" w14:textId="1AD7E468" w:rsidR="00076EA8" w:rsidRPr="00787C53" w:rsidRDefault="00252981" w:rsidP="00252981"><w:pPr><w:spacing w:after="0"/><w:jc w:val="right"/><w:rPr><w:noProof w:val="0"/><w:lang w:val="es-CO"/></w:rPr></w:pPr><w:r><w:rPr><w:noProof w:val="0"/><w:lang w:val="es-CO"/></w:rPr><w:t>Luisa Zambrano</w:t></w:r><w:r><w:rPr><w:noProof w:val="0"/><w:lang w:val="es-CO"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:noProof w:val="0"/><w:lang w:val="es-CO"/></w:rPr><w:t>- 201914911</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0"/><w:jc w:val="right"/><w:rPr><w:noProof w:val="0"/><w:lang w:val="es-CO"/></w:rPr></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Paragraphs(4).Range.InsertXML($xml11)

$xml12 = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="601EDB82" w14:textId="77777777" w:rsidR="00252981" w:rsidRDefault="00252981" w:rsidP="00252981"><w:pPr><w:spacing w:after="0"/><w:jc w:val="right"/><w:rPr><w:noProof w:val="0"/><w:lang w:val="es-CO"/></w:rPr></w:pPr><w:r><w:rPr><w:noProof w:val="0"/><w:lang w:val="es-CO"/></w:rPr><w:t>Samuel Jaramillo -</w:t></w:r><w:r><w:rPr><w:noProof w:val="0"/><w:lang w:val="es-CO"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:noProof w:val="0"/><w:lang w:val="es-CO"/></w:rPr><w:t>202010768</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Paragraphs(3).Range.InsertXML($xml12)

$xml13 = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="250F899C" w14:textId="77777777" w:rsidR="00252981" w:rsidRDefault="00252981" w:rsidP="00252981"><w:pPr><w:spacing w:after="0"/><w:jc w:val="right"/><w:rPr><w:noProof w:val="0"/><w:lang w:val="es-CO"/></w:rPr></w:pPr><w:r><w:rPr><w:noProof w:val="0"/><w:lang w:val="es-CO"/></w:rPr><w:t>Miguel Perdomo -</w:t></w:r><w:r><w:rPr><w:noProof w:val="0"/><w:lang w:val="es-CO"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:noProof w:val="0"/><w:lang w:val="es-CO"/></w:rPr><w:t>201913791</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Paragraphs(2).Range.InsertXML($xml13)

